# Add the "Abstract Class" vs "Virtual Class" comparison rows (29-30) to
# Sheet1, resize column B to fit the new (much longer) description text,
# and move the viewport/selection down to the newly-added rows — matching
# the author's "Add files via upload" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings are appended in this order so sharedStrings.xml ends
# up as: ... Abstract Class, Virtual Class , <abstract descr>, <virtual descr>
$ws.Range("A29").Value = "Abstract Class"
$ws.Range("A30").Value = "Virtual Class "
$ws.Range("B29").Value = "abstract class doesnt provide the implementations and forces derived class to override the method"
$ws.Range("B30").Value = "virtual class has implementation and provide derived class with option to override it"

# Column B needs to be much wider now that it holds long descriptions.
$ws.Columns.Item(2).ColumnWidth = 90

# Move the viewport down toward the new rows and park the selection on B33,
# same as the author's saved view state.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B33").Select()
